# NATMI output update ("Natmi following Dr Hou advice"):
# the Sending/Target cluster labels are corrected and the result table is
# expanded from a partial 6-row listing (rows 2-7) to the full 3x3
# Sending-cluster x Target-cluster matrix (ECs/FAPs/sCs) spanning rows 2-10,
# with recalculated statistics for every cluster pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs  Target=ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Clcf1"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.798783666666667
$ws.Range("H2").Value = 5.396351
$ws.Range("I2").Value = 0.2319744053785674
$ws.Range("J2").Value = 0.2319744053785674
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 37.42645866666667
$ws.Range("N2").Value = 112.279376
$ws.Range("O2").Value = 0.2415534622699011
$ws.Range("P2").Value = 0.2415534622699011
$ws.Range("Q2").Value = 67.32210255077511
$ws.Range("R2").Value = 605.898922956976
$ws.Range("S2").Value = 0.05603422077719453
$ws.Range("T2").Value = 0.05603422077719452

# Row 3: Sending=ECs  Target=FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Clcf1"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.798783666666667
$ws.Range("H3").Value = 5.396351
$ws.Range("I3").Value = 0.2319744053785674
$ws.Range("J3").Value = 0.2319744053785674
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 94.96115633333334
$ws.Range("N3").Value = 284.883469
$ws.Range("O3").Value = 0.6128871635375853
$ws.Range("P3").Value = 0.6128871635375853
$ws.Range("Q3").Value = 170.8145769801799
$ws.Range("R3").Value = 1537.331192821619
$ws.Range("S3").Value = 0.1421741353257882
$ws.Range("T3").Value = 0.1421741353257882

# Row 4: Sending=ECs  Target=sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Clcf1"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.798783666666667
$ws.Range("H4").Value = 5.396351
$ws.Range("I4").Value = 0.2319744053785674
$ws.Range("J4").Value = 0.2319744053785674
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.553069
$ws.Range("N4").Value = 67.659207
$ws.Range("O4").Value = 0.1455593741925136
$ws.Range("P4").Value = 0.1455593741925136
$ws.Range("Q4").Value = 40.56809215040633
$ws.Range("R4").Value = 365.112829353657
$ws.Range("S4").Value = 0.03376604927558472
$ws.Range("T4").Value = 0.03376604927558472

# Row 5: Sending=FAPs  Target=ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Clcf1"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.782700333333333
$ws.Range("H5").Value = 5.348101
$ws.Range("I5").Value = 0.2299002695301921
$ws.Range("J5").Value = 0.2299002695301921
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 37.42645866666667
$ws.Range("N5").Value = 112.279376
$ws.Range("O5").Value = 0.2415534622699011
$ws.Range("P5").Value = 0.2415534622699011
$ws.Range("Q5").Value = 66.7201603405529
$ws.Range("R5").Value = 600.481443064976
$ws.Range("S5").Value = 0.05553320608180135
$ws.Range("T5").Value = 0.05553320608180135

# Row 6: Sending=FAPs  Target=FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Clcf1"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.782700333333333
$ws.Range("H6").Value = 5.348101
$ws.Range("I6").Value = 0.2299002695301921
$ws.Range("J6").Value = 0.2299002695301921
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 94.96115633333334
$ws.Range("N6").Value = 284.883469
$ws.Range("O6").Value = 0.6128871635375853
$ws.Range("P6").Value = 0.6128871635375853
$ws.Range("Q6").Value = 169.2872850491521
$ws.Range("R6").Value = 1523.585565442369
$ws.Range("S6").Value = 0.1409029240888858
$ws.Range("T6").Value = 0.1409029240888858

# Row 7: Sending=FAPs  Target=sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Clcf1"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.782700333333333
$ws.Range("H7").Value = 5.348101
$ws.Range("I7").Value = 0.2299002695301921
$ws.Range("J7").Value = 0.2299002695301921
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.553069
$ws.Range("N7").Value = 67.659207
$ws.Range("O7").Value = 0.1455593741925136
$ws.Range("P7").Value = 0.1455593741925136
$ws.Range("Q7").Value = 40.20536362398966
$ws.Range("R7").Value = 361.848272615907
$ws.Range("S7").Value = 0.03346413935950496
$ws.Range("T7").Value = 0.03346413935950495

# Row 8: Sending=sCs  Target=ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Clcf1"
$ws.Range("C8").Value = "Il6st"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.172749333333334
$ws.Range("H8").Value = 12.518248
$ws.Range("I8").Value = 0.5381253250912406
$ws.Range("J8").Value = 0.5381253250912404
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 37.42645866666667
$ws.Range("N8").Value = 112.279376
$ws.Range("O8").Value = 0.2415534622699011
$ws.Range("P8").Value = 0.2415534622699011
$ws.Range("Q8").Value = 156.1712304503609
$ws.Range("R8").Value = 1405.541074053248
$ws.Range("S8").Value = 0.1299860354109052
$ws.Range("T8").Value = 0.1299860354109052

# Row 9: Sending=sCs  Target=FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Clcf1"
$ws.Range("C9").Value = "Il6st"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.172749333333334
$ws.Range("H9").Value = 12.518248
$ws.Range("I9").Value = 0.5381253250912406
$ws.Range("J9").Value = 0.5381253250912404
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 94.96115633333334
$ws.Range("N9").Value = 284.883469
$ws.Range("O9").Value = 0.6128871635375853
$ws.Range("P9").Value = 0.6128871635375853
$ws.Range("Q9").Value = 396.2491017824792
$ws.Range("R9").Value = 3566.241916042312
$ws.Range("S9").Value = 0.3298101041229114
$ws.Range("T9").Value = 0.3298101041229113

# Row 10: Sending=sCs  Target=sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Clcf1"
$ws.Range("C10").Value = "Il6st"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.172749333333334
$ws.Range("H10").Value = 12.518248
$ws.Range("I10").Value = 0.5381253250912406
$ws.Range("J10").Value = 0.5381253250912404
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.553069
$ws.Range("N10").Value = 67.659207
$ws.Range("O10").Value = 0.1455593741925136
$ws.Range("P10").Value = 0.1455593741925136
$ws.Range("Q10").Value = 94.10830363437067
$ws.Range("R10").Value = 846.9747327093361
$ws.Range("S10").Value = 0.07832918555742388
$ws.Range("T10").Value = 0.07832918555742387

